$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Mapping of row -> (new DAMSLTag, new DialogAct)
$changes = @{
    3   = @("sd", "Statement-non-opinion")
    6   = @("sd", "Statement-non-opinion")
    18  = @("sd", "Statement-non-opinion")
    37  = @("sd", "Statement-non-opinion")
    38  = @("aa", "Agree/Accept")
    43  = @("sd", "Statement-non-opinion")
    53  = @("sd", "Statement-non-opinion")
    57  = @("aa", "Agree/Accept")
    61  = @("sv", "Statement-opinion")
    71  = @("ba", "Appreciation")
    78  = @("sd", "Statement-non-opinion")
    86  = @("b", "Acknowledge (Backchannel)")
    114 = @("%", "Uninterpretable")
    132 = @("%", "Uninterpretable")
    137 = @("sd", "Statement-non-opinion")
    138 = @("%", "Uninterpretable")
    139 = @("sd", "Statement-non-opinion")
    141 = @("sd", "Statement-non-opinion")
    144 = @("ba", "Appreciation")
}

foreach ($row in $changes.Keys) {
    $vals = $changes[$row]
    $ws.Range("I$row").Value = $vals[0]
    $ws.Range("J$row").Value = $vals[1]
}
